$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

$ws.Range("J2").Value = "bron"
$ws.Range("J3").Value = "Bron"

For ($i = 4; $i -le 152; $i++) {
    $ws.Cells.Item($i, 10).Value = "VRS"
}

$ws.Range("J5:J152").Select()
